# Order History Session Storage
# Appends 4 new order history rows (10-13) to the "Order History" sheet,
# matching orders accumulated in session storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = 2; B = 45686.22928240741; C = "Karthik";       D = 2; E = 790; F = 10;  G = 5;   H = 15;  I = "Vanilla Shake (x4), Mango Lassi (x2)" },
    @{ A = 3; B = 45672.22928240741; C = "Reevan";        D = 4; E = 500; F = 3;   G = 4;   H = 7;   I = "Chicken Burger (x2), Chicken Cheese Burger (x2), Oreo Shake (x2), Vanilla Shake (x1), Mango Lassi (x4)" },
    @{ A = 4; B = 45658.22928240741; C = "Notsla Daniel"; D = 5; E = 345; F = 2.3; G = 1.2; H = 3.5; I = "Chicken Cheese Burger (x5)" },
    @{ A = 5; B = 45647.22928240741; C = "Karthik";       D = 6; E = 895; F = 3.9; G = 3.2; H = 7.1; I = "Chicken Cheese Pops (x3)" }
)

# Source cell whose date formatting (style) the new Date column should
# inherit - all existing date cells in column B already share this style.
$dateStyleSource = $ws.Cells.Item(2, 2)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A

    # Copy the date cell's formatting first so the new cell reuses the
    # existing date style instead of minting a brand-new one, then set
    # the actual value.
    $dateStyleSource.Copy($ws.Cells.Item($r, 2))
    $ws.Cells.Item($r, 2).Value = $row.B

    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
